# This script inserts two new weekly price rows for "Coliflor" (Cauliflower)
# right before the existing row 951, shifting all subsequent rows down by two.
# This mirrors the original data source being re-extracted with two additional
# weekly records (the most recent entries) inserted near the top of this
# particular date-block, while all previously existing rows keep their values
# and simply move down two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 951. Everything that used to be in
# rows 951.. moves down to 953..
$ws.Rows("951:952").Insert()

# ---- Row 951 (new) ----
$ws.Cells.Item(951, 1).Value2  = 6
$ws.Cells.Item(951, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(951, 3).Value2  = "Metropolitana"
$ws.Cells.Item(951, 4).Value2  = 44931
$ws.Cells.Item(951, 5).Value2  = 13
$ws.Cells.Item(951, 6).Value2  = 100112008
$ws.Cells.Item(951, 7).Value2  = "Coliflor"
$ws.Cells.Item(951, 8).Value2  = "Sin especificar"
$ws.Cells.Item(951, 9).Value2  = "Primera"
$ws.Cells.Item(951, 10).Value2 = 8200
$ws.Cells.Item(951, 11).Value2 = 700
$ws.Cells.Item(951, 12).Value2 = 800
$ws.Cells.Item(951, 13).Value2 = 740
$ws.Cells.Item(951, 14).Value2 = "`$/unidad"
$ws.Cells.Item(951, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(951, 16).Value2 = 740
$ws.Cells.Item(951, 17).Value2 = 1
$ws.Cells.Item(951, 18).Value2 = "Hortaliza"

# ---- Row 952 (new) ----
$ws.Cells.Item(952, 1).Value2  = 6
$ws.Cells.Item(952, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(952, 3).Value2  = "Metropolitana"
$ws.Cells.Item(952, 4).Value2  = 44931
$ws.Cells.Item(952, 5).Value2  = 13
$ws.Cells.Item(952, 6).Value2  = 100112008
$ws.Cells.Item(952, 7).Value2  = "Coliflor"
$ws.Cells.Item(952, 8).Value2  = "Sin especificar"
$ws.Cells.Item(952, 9).Value2  = "Segunda"
$ws.Cells.Item(952, 10).Value2 = 2800
$ws.Cells.Item(952, 11).Value2 = 500
$ws.Cells.Item(952, 12).Value2 = 500
$ws.Cells.Item(952, 13).Value2 = 500
$ws.Cells.Item(952, 14).Value2 = "`$/unidad"
$ws.Cells.Item(952, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(952, 16).Value2 = 500
$ws.Cells.Item(952, 17).Value2 = 1
$ws.Cells.Item(952, 18).Value2 = "Hortaliza"

# Make sure column D (Fecha) keeps the date style used elsewhere in that
# column (style carries over from Insert, but set explicitly to be safe).
$ws.Range("D951:D952").NumberFormat = $ws.Range("D950").NumberFormat
